$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '244.12'

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '23.81'

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.05812'

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '6.468'

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.226'

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.8080'

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.8812'

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.1392'

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07093'

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.03194'

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.03047'

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.09332'

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '3.826'

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.001536'

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.04706'

# Row 18
$ws.Cells.Item(18, 2).Value = 'One'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.0006023'
$ws.Cells.Item(18, 5).Value = '17OneONE'

# Row 19
$ws.Cells.Item(19, 2).Value = 'TigerCash'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.006162'
$ws.Cells.Item(19, 5).Value = '18TigerCashTCH'

# Row 20
$ws.Cells.Item(20, 2).Value = 'BitKan'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.001257'
$ws.Cells.Item(20, 5).Value = '19BitKanKAN'

# Row 21
$ws.Cells.Item(21, 2).Value = 'HotbitToken'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.004069'
$ws.Cells.Item(21, 5).Value = '20HotbitTokenHTB'

# Row 22
$ws.Cells.Item(22, 2).Value = 'NitroEx'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.00008702'
$ws.Cells.Item(22, 5).Value = '21NitroExNTX'

# Row 23
$ws.Cells.Item(23, 2).Value = 'LEO'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '3.545'
$ws.Cells.Item(23, 5).Value = '22LEOLEO'

# Row 24
$ws.Cells.Item(24, 2).Value = 'BTSEToken'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.155'
$ws.Cells.Item(24, 5).Value = '23BTSETokenBTSE'

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.3186'

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.1328'

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.0002330'

# Row 41
$ws.Cells.Item(41, 2).Value = 'BKEXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.1051'
$ws.Cells.Item(41, 5).Value = '40BKEXTokenBKK'

# Row 42
$ws.Cells.Item(42, 2).Value = 'CEJI'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.002380'
$ws.Cells.Item(42, 5).Value = '41CEJICEJI'

# Row 43
$ws.Cells.Item(43, 2).Value = 'KickToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.003246'
$ws.Cells.Item(43, 5).Value = '42KickTokenKICKWorstin24h'

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.007869'

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.00005320'

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.5353'
$ws.Cells.Item(47, 5).Value = '46CoinbaseStockTokenCOIN'

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.002592'
$ws.Cells.Item(48, 5).Value = '47BOLOBOLOBestin24h'
